$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: D2872_2_bg_detlim
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("D2872_2_bg_detlim")

# Rows 2 and 5 share identical values
foreach ($r in 2,5) {
    $ws1.Range("B$r").Value = 0.013
    $ws1.Range("C$r").Value = 0.011
    $ws1.Range("D$r").Value = 0.012
    $ws1.Range("E$r").Value = 0.029
    $ws1.Range("F$r").Value = 0.015
    $ws1.Range("G$r").Value = 0.016
    $ws1.Range("H$r").Value = 0.014
    $ws1.Range("I$r").Value = 0.014
    $ws1.Range("J$r").Value = 0.015
    $ws1.Range("K$r").Value = 0.033
    $ws1.Range("L$r").Value = 0.026
    $ws1.Range("M$r").Value = 0.03
    $ws1.Range("N$r").Value = 0.028
    $ws1.Range("O$r").Value = 0.03
    $ws1.Range("P$r").Value = 0.02
    $ws1.Range("Q$r").Value = 0.008
    $ws1.Range("R$r").Value = 0.011
    $ws1.Range("S$r").Value = 0.033
}

# Row 8
$ws1.Range("B8").Value = 0.017
$ws1.Range("C8").Value = 0.015
$ws1.Range("D8").Value = 0.016
$ws1.Range("E8").Value = 0.039
$ws1.Range("F8").Value = 0.02
$ws1.Range("G8").Value = 0.021
$ws1.Range("H8").Value = 0.02
$ws1.Range("I8").Value = 0.02
$ws1.Range("J8").Value = 0.02
$ws1.Range("K8").Value = 0.044
$ws1.Range("L8").Value = 0.035
$ws1.Range("M8").Value = 0.041
$ws1.Range("N8").Value = 0.038
$ws1.Range("O8").Value = 0.041
$ws1.Range("P8").Value = 0.028
$ws1.Range("Q8").Value = 0.011
$ws1.Range("R8").Value = 0.015
$ws1.Range("S8").Value = 0.044

# ------------------------------------------------------------------
# Sheet 2: D2872_3_bg_apf_detlim
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("D2872_3_bg_apf_detlim")

# Rows 2 and 5 share identical values
foreach ($r in 2,5) {
    $ws2.Range("B$r").Value = 0.015
    $ws2.Range("C$r").Value = 0.013
    $ws2.Range("D$r").Value = 0.015
    $ws2.Range("E$r").Value = 0.035
    $ws2.Range("F$r").Value = 0.018
    $ws2.Range("G$r").Value = 0.019
    $ws2.Range("H$r").Value = 0.017
    $ws2.Range("I$r").Value = 0.017
    $ws2.Range("J$r").Value = 0.018
    $ws2.Range("K$r").Value = 0.039
    $ws2.Range("L$r").Value = 0.031
    $ws2.Range("M$r").Value = 0.037
    $ws2.Range("N$r").Value = 0.034
    $ws2.Range("O$r").Value = 0.036
    $ws2.Range("P$r").Value = 0.025
    $ws2.Range("Q$r").Value = 0.01
    $ws2.Range("R$r").Value = 0.013
    $ws2.Range("S$r").Value = 0.039
}

# Row 8
$ws2.Range("B8").Value = 0.021
$ws2.Range("C8").Value = 0.018
$ws2.Range("D8").Value = 0.02
$ws2.Range("E8").Value = 0.047
$ws2.Range("F8").Value = 0.024
$ws2.Range("G8").Value = 0.025
$ws2.Range("H8").Value = 0.024
$ws2.Range("I8").Value = 0.024
$ws2.Range("J8").Value = 0.024
$ws2.Range("K8").Value = 0.053
$ws2.Range("L8").Value = 0.042
$ws2.Range("M8").Value = 0.049
$ws2.Range("N8").Value = 0.045
$ws2.Range("O8").Value = 0.049
$ws2.Range("P8").Value = 0.033
$ws2.Range("Q8").Value = 0.013
$ws2.Range("R8").Value = 0.018
$ws2.Range("S8").Value = 0.053
